# Scaffold de algunas paginas.
# Applies the worksheet-content portion of the authored change:
#  - a thin column of empty, center-aligned cells (G15:G23) under the
#    "Proyectos" mini-table
#  - a new "Relacion Tecnico/Rol" scaffold table in B41:D49 with a helper
#    column (F) that builds the literal C# initializer text for each row
#  - updates the current selection to match the author's final cursor spot

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) G15:G23 -> empty cells, center/center aligned (reuses existing style)
# ---------------------------------------------------------------------
$gCol = $ws.Range("G15:G23")
$gCol.HorizontalAlignment = -4108   # xlCenter
$gCol.VerticalAlignment = -4108     # xlCenter

# ---------------------------------------------------------------------
# 2) New scaffold table, B41:D49 (+ header) — copy borders/alignment from
#    the look-alike table that already exists at row 30 so no new border
#    definitions get minted, only the formats we actually need.
# ---------------------------------------------------------------------
$ws.Range("B31").Copy() | Out-Null
$ws.Range("B41:D41").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B5").Copy() | Out-Null
$ws.Range("B42:B49").PasteSpecial(-4122) | Out-Null
$ws.Range("D42:D49").PasteSpecial(-4122) | Out-Null

$ws.Range("B31").Copy() | Out-Null
$ws.Range("C42:C49").PasteSpecial(-4122) | Out-Null

$ws.Range("B41").Value = "TecnicoId"
$ws.Range("C41").Value = "SolicitudId"
$ws.Range("D41").Value = "RolRequerido"

$tecnicos = @(
  "57bf6b3f-26f0-4eaa-9f66-14b3e6fdfce2",
  "57bf6b3f-26f0-4eaa-9f66-14b3e6fdfce2",
  "0626bd2e-c394-4f89-bb52-8dcf01b0128c",
  "0626bd2e-c394-4f89-bb52-8dcf01b0128c",
  "0626bd2e-c394-4f89-bb52-8dcf01b0128c",
  "cf374546-893e-4b69-8622-a334fb02ade8",
  "cf374546-893e-4b69-8622-a334fb02ade8",
  "cf374546-893e-4b69-8622-a334fb02ade8"
)
$roles = @(
  "Operador de Cabina 03 y Estudio de Radio",
  "Sonidista",
  "Presentador / conductor",
  "Sonidista",
  "Redactor creativo",
  "Operador de Cabina 02",
  "Diseñador gráfico",
  "Cámara y asistente de cámara"
)

for ($i = 0; $i -lt 8; $i++) {
    $row = 42 + $i
    $ws.Cells.Item($row, 2).Value = $tecnicos[$i]
    $ws.Cells.Item($row, 3).Value = $i + 1
    $ws.Cells.Item($row, 4).Value = $roles[$i]
}

# Helper column F: literal C# object-initializer string, left/center
# aligned (a brand-new style — no existing cell already has this combo).
$fCol = $ws.Range("F42:F49")
$fCol.HorizontalAlignment = -4131   # xlLeft
$fCol.VerticalAlignment = -4108     # xlCenter

$ws.Range("F42").Formula = '="new RelacionTecnicoRol {TecnicoId = " & B42 & ", RolId = " & D42 & "},"'
$ws.Range("F43:F49").Formula = '="new RelacionTecnicoRol {TecnicoId = " & B43 & ", RolId = " & D43 & "},"'

# ---------------------------------------------------------------------
# 3) Final cursor position, matching the author's saved view.
# ---------------------------------------------------------------------
$ws.Range("C60").Select() | Out-Null
